$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the row 2 / row 3 values for columns D, M, N, O, P, S
# New row2: date 44421 (was row3's date), M=30, N=24000, O=24000, P=24000, S=1200
# New row3: date 44291 (was row2's date), M=15, N=23000, O=23000, P=23000, S=1150

$ws.Range("D2").Value = 44421
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("S2").Value = 1200

$ws.Range("D3").Value = 44291
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 23000
$ws.Range("P3").Value = 23000
$ws.Range("S3").Value = 1150
